$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are a contract negotiator for a fabric supplier. A furniture manufacturing company has expressed interest in purchasing custom-made fabric of a particular quality from the supplier. The standard pricing for the fabric is negotiable, but the proposed price by the company is significantly lower than what you would typically charge for fabric of that quality. Among the following options, what would be the best response?",
        "ques_type": 2,
        "options": [
            "Suggest a counteroffer that is closer to the standard price for the fabric.",
            "Accept the price offered by the client and proceed to make the fabric.",
            "Reject the price offered by the client and ask him to reconsider.",
            "Disregard the offer and make the fabric according to the standard price."
        ],
        "score": "Suggest a counteroffer that is closer to the standard price for the fabric."
    },
    {
        "title": "You own a manufacturing business and order 100 units of a specific raw material. The contract stipulates the quality and grade of the material to be delivered. Upon delivery, you notice that the material is of lower quality and grade than specified. You inform the seller of the non-compliance and request that they remedy the situation. However, the seller refuses to take any action. What can you do to recover your losses?",
        "ques_type": 2,
        "options": [
            "Sue the seller for specific performance of the contract, forcing them to provide the goods that were promised in the contract.",
            "Waive the breach and continue with the contract, asking the seller to fulfill their obligations in the future.",
            "Terminate the contract and sue for compensatory damages for the losses incurred due to the non-compliance of the goods.",
            "Sue the seller for punitive damages to punish them for the breach of contract and to deter future breaches."
        ],
        "score": "Terminate the contract and sue for compensatory damages for the losses incurred due to the non-compliance of the goods."
    },
    {
        "title": "You are a project manager for a construction company that signed a contract with a vendor to provide a specialized crane to complete a high-rise building project. However, the vendor notifies you that they can no longer deliver the crane, as their factory was destroyed by a natural disaster. What is the most appropriate course of action?",
        "ques_type": 2,
        "options": [
            "Negotiate with the vendor for a later delivery date. ",
            "Continue working with the vendor and ask for an alternative piece of equipment. ",
            "Waive the breach and continue with the contract.",
            "Terminate the contract and seek a new vendor. "
        ],
        "score": "Terminate the contract and seek a new vendor."
    },
    {
        "title": "As a procurement manager, you're finalizing a contract with a vendor. The product they offer has an industry-standard definition. However, you've used a term in the contract that deviates from this common meaning. The vendor has signed the contract without raising any objections.How should you interpret the term that differs from its commonly understood meaning?",
        "ques_type": 2,
        "options": [
            "Assume your definition is accepted because no objections were raised.",
            "Consult with legal experts or industry professionals to determine the intended meaning of the term and ensure adherence to the industry-standard definition.",
            "Stick to the industry-standard definition despite your different usage in the contract.",
            "Seek clarification from the vendor to avoid misunderstandings and potential legal disputes in the future."
        ],
        "score": "Seek clarification from the vendor to avoid misunderstandings and potential legal disputes in the future."
    }
]
'@

$ws.Cells.Item(2,1).EntireRow.Delete() | Out-Null
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).EntireRow.AutoFit()
